$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DAMSLTag (column I) and DialogAct (column J) values for rows
# re-annotated by re-running SGNN after transcript clean-up.
$ws.Range("I6").Value = 'sv'
$ws.Range("J6").Value = 'Statement-opinion'
$ws.Range("I20").Value = 'sv'
$ws.Range("J20").Value = 'Statement-opinion'
$ws.Range("I24").Value = 'sv'
$ws.Range("J24").Value = 'Statement-opinion'
$ws.Range("I30").Value = 'sv'
$ws.Range("J30").Value = 'Statement-opinion'
$ws.Range("I33").Value = 'sv'
$ws.Range("J33").Value = 'Statement-opinion'
$ws.Range("I37").Value = 'ba'
$ws.Range("J37").Value = 'Appreciation'
$ws.Range("I39").Value = 'aa'
$ws.Range("J39").Value = 'Agree/Accept'
$ws.Range("I66").Value = 'sd'
$ws.Range("J66").Value = 'Statement-non-opinion'
$ws.Range("I78").Value = 'b'
$ws.Range("J78").Value = 'Acknowledge (Backchannel)'
$ws.Range("I80").Value = 'sd'
$ws.Range("J80").Value = 'Statement-non-opinion'
$ws.Range("I82").Value = 'sd'
$ws.Range("J82").Value = 'Statement-non-opinion'
$ws.Range("I83").Value = 'aa'
$ws.Range("J83").Value = 'Agree/Accept'
$ws.Range("I85").Value = 'b'
$ws.Range("J85").Value = 'Acknowledge (Backchannel)'
$ws.Range("I96").Value = '%'
$ws.Range("J96").Value = 'Uninterpretable'
$ws.Range("I105").Value = 'sv'
$ws.Range("J105").Value = 'Statement-opinion'
$ws.Range("I107").Value = 'sv'
$ws.Range("J107").Value = 'Statement-opinion'
$ws.Range("I110").Value = 'sv'
$ws.Range("J110").Value = 'Statement-opinion'
$ws.Range("I118").Value = 'sd'
$ws.Range("J118").Value = 'Statement-non-opinion'
$ws.Range("I122").Value = 'sd'
$ws.Range("J122").Value = 'Statement-non-opinion'
$ws.Range("I123").Value = 'sv'
$ws.Range("J123").Value = 'Statement-opinion'
$ws.Range("I126").Value = 'sv'
$ws.Range("J126").Value = 'Statement-opinion'
$ws.Range("I138").Value = 'sv'
$ws.Range("J138").Value = 'Statement-opinion'
$ws.Range("I182").Value = 'sv'
$ws.Range("J182").Value = 'Statement-opinion'
$ws.Range("I189").Value = 'sd'
$ws.Range("J189").Value = 'Statement-non-opinion'
$ws.Range("I207").Value = 'sd'
$ws.Range("J207").Value = 'Statement-non-opinion'
$ws.Range("I216").Value = 'aa'
$ws.Range("J216").Value = 'Agree/Accept'
$ws.Range("I217").Value = 'sd'
$ws.Range("J217").Value = 'Statement-non-opinion'
$ws.Range("I224").Value = 'sv'
$ws.Range("J224").Value = 'Statement-opinion'
$ws.Range("I229").Value = '%'
$ws.Range("J229").Value = 'Uninterpretable'
$ws.Range("I231").Value = '%'
$ws.Range("J231").Value = 'Uninterpretable'
$ws.Range("I242").Value = 'b'
$ws.Range("J242").Value = 'Acknowledge (Backchannel)'
$ws.Range("I246").Value = 'sv'
$ws.Range("J246").Value = 'Statement-opinion'
$ws.Range("I254").Value = 'b'
$ws.Range("J254").Value = 'Acknowledge (Backchannel)'
$ws.Range("I259").Value = 'aa'
$ws.Range("J259").Value = 'Agree/Accept'
$ws.Range("I274").Value = 'b'
$ws.Range("J274").Value = 'Acknowledge (Backchannel)'
$ws.Range("I276").Value = 'sv'
$ws.Range("J276").Value = 'Statement-opinion'
$ws.Range("I280").Value = 'sv'
$ws.Range("J280").Value = 'Statement-opinion'
